$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = "Gestão"
$ws.Range("E3").Value = "Desenho Técnico"

# Row 4
$ws.Range("C4").Value = "-"
$ws.Range("E4").Value = "Tecnologia dos Materiais"
$ws.Range("F4").Value = "Tecnologia dos Materiais"

# Row 6
$ws.Range("C6").Value = "Desenho Técnico"
$ws.Range("D6").Value = "Circuitos Elétricos"
